$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the "date group" data between rows 2-5 (La Ligua, 2022-11-25)
# and rows 6-9 (Provincia de Limari, 2022-12-13). Columns D, M, N, O, P, R, S
# are swapped pairwise: row2<->row6, row3<->row7, row4<->row8, row5<->row9.

$pairs = @(2, 6), @(3, 7), @(4, 8), @(5, 9)
$cols = 4, 13, 14, 15, 16, 18, 19

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2

        $ws.Cells.Item($r1, $col).Value2 = $v2
        $ws.Cells.Item($r2, $col).Value2 = $v1
    }
}
